$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 223-255 (33 rows of invalid/offshore territory country codes
# that are no longer part of the controlled vocabulary mapping).
$ws.Range("A223:D255").EntireRow.Delete()

# Restore the view state (scroll position/selection) to reflect where the
# author was working after the cleanup.
$ws.Range("H220").Select()
$excel.ActiveWindow.ScrollRow = 208
$excel.ActiveWindow.ScrollColumn = 1
